$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and report week dates) ---
$ws.Range("A8").Characters(21, 2).Text = "33"
$ws.Range("C9").Characters(27, 8).Text = "8/12/2024"
$ws.Range("C9").Characters(47, 9).Text = "8/18/2024"

# --- Column E width change ---
$ws.Columns.Item(5).ColumnWidth = 6.168446

# --- Type-changing cells: copy style from a same-typed donor cell, then set value ---
# C22: text("0") -> number
$ws.Range("D22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value2 = 4

# C28: text("0") -> number
$ws.Range("F28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value2 = 1

# D28: number -> text("0")
$ws.Range("D33").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value2 = "0"

# E28: number -> text("***.*")
$ws.Range("E33").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value2 = "***.*"

# C33: number -> text("0")
$ws.Range("D33").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("C33").Value2 = "0"

# --- Plain numeric value updates ---
# Row 15
$ws.Range("C15").Value2 = 1
$ws.Range("F15").Value2 = 4
$ws.Range("I15").Value2 = 11
$ws.Range("K15").Value2 = -38.888888888888
$ws.Range("L15").Value2 = -26.666666666666
$ws.Range("M15").Value2 = 10
$ws.Range("N15").Value2 = 57.142857142857
# Row 16
$ws.Range("C16").Value2 = 6
$ws.Range("D16").Value2 = 7
$ws.Range("E16").Value2 = -14.285714285714
$ws.Range("F16").Value2 = 22
$ws.Range("G16").Value2 = 19
$ws.Range("H16").Value2 = 15.78947368421
$ws.Range("I16").Value2 = 160
$ws.Range("J16").Value2 = 168
$ws.Range("K16").Value2 = -4.761904761904
$ws.Range("L16").Value2 = 25
$ws.Range("M16").Value2 = 44.144144144144
$ws.Range("N16").Value2 = -75.720789074355
# Row 17
$ws.Range("C17").Value2 = 7
$ws.Range("D17").Value2 = 2
$ws.Range("E17").Value2 = 250
$ws.Range("F17").Value2 = 22
$ws.Range("G17").Value2 = 15
$ws.Range("H17").Value2 = 46.666666666666
$ws.Range("I17").Value2 = 162
$ws.Range("J17").Value2 = 147
$ws.Range("K17").Value2 = 10.204081632653
$ws.Range("L17").Value2 = 11.724137931034
$ws.Range("M17").Value2 = 102.5
$ws.Range("N17").Value2 = -15.183246073298
# Row 18
$ws.Range("C18").Value2 = 8
$ws.Range("D18").Value2 = 8
$ws.Range("E18").Value2 = 0
$ws.Range("F18").Value2 = 25
$ws.Range("G18").Value2 = 21
$ws.Range("H18").Value2 = 19.047619047619
$ws.Range("I18").Value2 = 166
$ws.Range("J18").Value2 = 128
$ws.Range("K18").Value2 = 29.6875
$ws.Range("L18").Value2 = 48.214285714285
$ws.Range("M18").Value2 = 7.096774193548
$ws.Range("N18").Value2 = -81.995661605206
# Row 19
$ws.Range("D19").Value2 = 15
$ws.Range("E19").Value2 = 6.666666666666
$ws.Range("F19").Value2 = 56
$ws.Range("G19").Value2 = 49
$ws.Range("H19").Value2 = 14.285714285714
$ws.Range("I19").Value2 = 467
$ws.Range("J19").Value2 = 471
$ws.Range("K19").Value2 = -0.849256900212
$ws.Range("L19").Value2 = 9.624413145539
$ws.Range("M19").Value2 = 62.152777777777
$ws.Range("N19").Value2 = -17.781690140845
# Row 20
$ws.Range("D20").Value2 = 11
$ws.Range("E20").Value2 = -54.545454545454
$ws.Range("F20").Value2 = 33
$ws.Range("G20").Value2 = 33
$ws.Range("H20").Value2 = 0
$ws.Range("I20").Value2 = 160
$ws.Range("J20").Value2 = 165
$ws.Range("K20").Value2 = -3.030303030303
$ws.Range("L20").Value2 = 34.453781512605
$ws.Range("M20").Value2 = 16.788321167883
$ws.Range("N20").Value2 = -87.291501191421
# Row 21
$ws.Range("C21").Value2 = 43
$ws.Range("D21").Value2 = 43
$ws.Range("E21").Value2 = 0
$ws.Range("F21").Value2 = 162
$ws.Range("G21").Value2 = 137
$ws.Range("H21").Value2 = 18.248175182481
$ws.Range("I21").Value2 = 1127
$ws.Range("J21").Value2 = 1097
$ws.Range("K21").Value2 = 2.734731084776
$ws.Range("L21").Value2 = 19.259259259259
$ws.Range("M21").Value2 = 43.933588761175
$ws.Range("N21").Value2 = -68.815716657443
# Row 22
$ws.Range("E22").Value2 = 300
$ws.Range("F22").Value2 = 6
$ws.Range("G22").Value2 = 4
$ws.Range("H22").Value2 = 50
$ws.Range("I22").Value2 = 38
$ws.Range("J22").Value2 = 56
$ws.Range("K22").Value2 = -32.142857142857
$ws.Range("L22").Value2 = -19.148936170212
$ws.Range("M22").Value2 = 35.714285714285
# Row 24
$ws.Range("D24").Value2 = 45
$ws.Range("E24").Value2 = 8.888888888888
$ws.Range("F24").Value2 = 186
$ws.Range("G24").Value2 = 227
$ws.Range("H24").Value2 = -18.06167400881
$ws.Range("I24").Value2 = 1368
$ws.Range("J24").Value2 = 1251
$ws.Range("K24").Value2 = 9.352517985611
$ws.Range("L24").Value2 = 55.454545454545
$ws.Range("M24").Value2 = 129.145728643216
# Row 25
$ws.Range("C25").Value2 = 28
$ws.Range("D25").Value2 = 33
$ws.Range("E25").Value2 = -15.151515151515
$ws.Range("F25").Value2 = 118
$ws.Range("G25").Value2 = 153
$ws.Range("H25").Value2 = -22.875816993464
$ws.Range("I25").Value2 = 926
$ws.Range("J25").Value2 = 793
$ws.Range("K25").Value2 = 16.771752837326
$ws.Range("L25").Value2 = 174.777448071217
# Row 26
$ws.Range("C26").Value2 = 6
$ws.Range("D26").Value2 = 8
$ws.Range("E26").Value2 = -25
$ws.Range("F26").Value2 = 30
$ws.Range("G26").Value2 = 32
$ws.Range("H26").Value2 = -6.25
$ws.Range("I26").Value2 = 369
$ws.Range("J26").Value2 = 319
$ws.Range("K26").Value2 = 15.673981191222
$ws.Range("L26").Value2 = 0.271739130434
$ws.Range("M26").Value2 = 16.037735849056
# Row 27
$ws.Range("C27").Value2 = 1
$ws.Range("F27").Value2 = 5
$ws.Range("I27").Value2 = 18
$ws.Range("K27").Value2 = -21.739130434782
$ws.Range("L27").Value2 = -10
# Row 28
$ws.Range("F28").Value2 = 3
$ws.Range("H28").Value2 = -50
# Row 33
